$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 19
$ws.Range("H19").Value = 398.5
$ws.Range("J19").Value = 297.2
$ws.Range("L19").Value = 297.2
$ws.Range("N19").Value = -647.2
# Row 54
$ws.Range("H54").Value = 0
$ws.Range("J54").Value = 0
$ws.Range("L54").Value = 0
$ws.Range("N54").ClearContents()
# Row 74
$ws.Range("H74").Value = 5999
$ws.Range("I74").Value = 5999
$ws.Range("K74").Value = 5999
$ws.Range("M74").Value = -5063
# Row 77
$ws.Range("H77").Value = 5999
$ws.Range("I77").Value = 5999
$ws.Range("K77").Value = 29995
$ws.Range("M77").Value = -25315
# Row 121
$ws.Range("H121").Value = 631.62964
$ws.Range("I121").Value = 230
$ws.Range("J121").Value = 647.0769
$ws.Range("K121").Value = 690
$ws.Range("L121").Value = 1941.2307
$ws.Range("M121").Value = 1057
$ws.Range("N121").Value = -5435.2307
# Row 137
$ws.Range("H137").Value = 2088.889
$ws.Range("I137").Value = 1912.9166
$ws.Range("J137").Value = 2440.8333
$ws.Range("K137").Value = 5738.7498
$ws.Range("L137").Value = 7322.499899999999
$ws.Range("M137").Value = -3188.7498
$ws.Range("N137").Value = -12422.4999
# Row 138
$ws.Range("H138").Value = 2771.8489
$ws.Range("I138").Value = 1790.091
$ws.Range("J138").Value = 2915.84
$ws.Range("K138").Value = 5370.272999999999
$ws.Range("L138").Value = 8747.52
$ws.Range("M138").Value = -230.2729999999992
$ws.Range("N138").Value = -19027.52

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 74
$ws.Range("H74").Value = 1907.76
$ws.Range("I74").Value = 1067.5333
$ws.Range("K74").Value = 1067.5333
$ws.Range("M74").Value = -193.5333000000001
# Row 77
$ws.Range("H77").Value = 1907.76
$ws.Range("I77").Value = 1067.5333
$ws.Range("K77").Value = 5337.6665
$ws.Range("M77").Value = -969.6665000000003
# Row 97
$ws.Range("H97").Value = 6206.222
$ws.Range("I97").Value = 682.125
$ws.Range("J97").Value = 50399
$ws.Range("K97").Value = 682.125
$ws.Range("L97").Value = 50399
$ws.Range("M97").Value = -186.125
$ws.Range("N97").Value = -51391
# Row 102
$ws.Range("H102").Value = 16670927
$ws.Range("I102").Value = 23812182
$ws.Range("K102").Value = 23812182
$ws.Range("M102").Value = -23810560
# Row 122
$ws.Range("H122").Value = 2170.1304
$ws.Range("I122").Value = 1810.3334
$ws.Range("J122").Value = 3465.4
$ws.Range("K122").Value = 5431.0002
$ws.Range("L122").Value = 10396.2
$ws.Range("M122").Value = -2981.0002
$ws.Range("N122").Value = -15296.2
# Row 132
$ws.Range("H132").Value = 2852.3794
$ws.Range("I132").Value = 2686.7144
$ws.Range("J132").Value = 3287.25
$ws.Range("K132").Value = 8060.1432
$ws.Range("L132").Value = 9861.75
$ws.Range("M132").Value = -5530.1432
$ws.Range("N132").Value = -14921.75

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 58
$ws.Range("H58").Value = 28348.8
$ws.Range("I58").Value = 0
$ws.Range("J58").Value = 28348.8
$ws.Range("K58").Value = 0
$ws.Range("L58").Value = 28348.8
$ws.Range("M58").ClearContents()
$ws.Range("N58").Value = -28936.8
# Row 94
$ws.Range("H94").Value = 7143578
$ws.Range("I94").Value = 8621364
$ws.Range("J94").Value = 945
$ws.Range("K94").Value = 8621364
$ws.Range("L94").Value = 945
$ws.Range("M94").Value = -8620913
$ws.Range("N94").Value = -1847
# Row 129
$ws.Range("H129").Value = 49999
$ws.Range("J129").Value = 49999
$ws.Range("L129").Value = 49999
$ws.Range("N129").Value = -59999
# Row 134
$ws.Range("H134").Value = 9817.538
$ws.Range("I134").Value = 1233.1428
$ws.Range("J134").Value = 19832.666
$ws.Range("K134").Value = 3699.4284
$ws.Range("L134").Value = 59497.99800000001
$ws.Range("M134").Value = -1164.4284
$ws.Range("N134").Value = -64567.99800000001

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 2154.1943
$ws.Range("I31").Value = 2072.9143
$ws.Range("J31").Value = 4999
$ws.Range("K31").Value = 2072.9143
$ws.Range("L31").Value = 4999
$ws.Range("M31").Value = -1777.9143
$ws.Range("N31").Value = -5589
# Row 34
$ws.Range("H34").Value = 2154.1943
$ws.Range("I34").Value = 2072.9143
$ws.Range("J34").Value = 4999
$ws.Range("K34").Value = 2072.9143
$ws.Range("L34").Value = 4999
$ws.Range("M34").Value = -1870.9143
$ws.Range("N34").Value = -5403
# Row 132
$ws.Range("H132").Value = 2316.875
$ws.Range("I132").Value = 2001.3043
$ws.Range("J132").Value = 3123.3333
$ws.Range("K132").Value = 6003.9129
$ws.Range("L132").Value = 9369.999899999999
$ws.Range("M132").Value = -3473.9129
$ws.Range("N132").Value = -14429.9999
# Row 134
$ws.Range("H134").Value = 14287412
$ws.Range("I134").Value = 1775.3636
$ws.Range("J134").Value = 38463104
$ws.Range("K134").Value = 5326.0908
$ws.Range("L134").Value = 115389312
$ws.Range("M134").Value = -2791.0908
$ws.Range("N134").Value = -115394382

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 5
$ws.Range("H5").Value = 1631.6364
$ws.Range("I5").Value = 1538.9445
$ws.Range("J5").Value = 2048.75
$ws.Range("K5").Value = 4616.833500000001
$ws.Range("L5").Value = 6146.25
$ws.Range("M5").Value = -4504.833500000001
$ws.Range("N5").Value = -6370.25
# Row 42
$ws.Range("H42").Value = 3788.5715
$ws.Range("J42").Value = 3788.5715
$ws.Range("L42").Value = 11365.7145
$ws.Range("N42").Value = -12433.7145
# Row 122
$ws.Range("H122").Value = 777.7381
$ws.Range("J122").Value = 819.29034
$ws.Range("L122").Value = 7373.61306
$ws.Range("N122").Value = -12273.61306
# Row 135
$ws.Range("H135").Value = 1631.6364
$ws.Range("I135").Value = 1538.9445
$ws.Range("J135").Value = 2048.75
$ws.Range("K135").Value = 13850.5005
$ws.Range("L135").Value = 18438.75
$ws.Range("M135").Value = -11315.5005
$ws.Range("N135").Value = -23508.75

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 126
$ws.Range("H126").Value = 2006.1333
$ws.Range("I126").Value = 1779.1111
$ws.Range("J126").Value = 2346.6667
$ws.Range("K126").Value = 5337.3333
$ws.Range("L126").Value = 7040.000100000001
$ws.Range("M126").Value = -2867.3333
$ws.Range("N126").Value = -11980.0001

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 46
$ws.Range("H46").Value = 5931.875
$ws.Range("I46").Value = 730
$ws.Range("J46").Value = 7665.8335
$ws.Range("K46").Value = 730
$ws.Range("L46").Value = 7665.8335
$ws.Range("M46").Value = -542
$ws.Range("N46").Value = -8041.8335

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 132
$ws.Range("H132").Value = 4061.1082
$ws.Range("I132").Value = 4353.931
$ws.Range("J132").Value = 2999.625
$ws.Range("K132").Value = 13061.793
$ws.Range("L132").Value = 8998.875
$ws.Range("M132").Value = -10531.793
$ws.Range("N132").Value = -14058.875
